$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("H2").Value = 7243
$ws.Range("I2").Value = 3490
$ws.Range("I3").Value = 3619
$ws.Range("G4").Value = 1435
$ws.Range("H4").Value = 1667
$ws.Range("I4").Value = 844
$ws.Range("I5").Value = 337
$ws.Range("I6").Value = 4057
$ws.Range("G7").Value = 24657
$ws.Range("H7").Value = 25978
$ws.Range("I7").Value = 12347

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I6").Value = 85
$ws.Range("G7").Value = 708
$ws.Range("I8").Value = 757
$ws.Range("I12").Value = 24
$ws.Range("I15").Value = 146
$ws.Range("I18").Value = 82
$ws.Range("I19").Value = 324
$ws.Range("I22").Value = 34
$ws.Range("I26").Value = 15
$ws.Range("I27").Value = 110
$ws.Range("I29").Value = 802
$ws.Range("I33").Value = 557
$ws.Range("I37").Value = 391
$ws.Range("I42").Value = 426
$ws.Range("I44").Value = 92
$ws.Range("I46").Value = 28
$ws.Range("I50").Value = 56
$ws.Range("I51").Value = 120
$ws.Range("I52").Value = 270
$ws.Range("I53").Value = 138
$ws.Range("I54").Value = 276
$ws.Range("I55").Value = 140
$ws.Range("H63").Value = 207
$ws.Range("I63").Value = 47
$ws.Range("I65").Value = 273
$ws.Range("I67").Value = 476
$ws.Range("I75").Value = 43
$ws.Range("I76").Value = 188
$ws.Range("I78").Value = 174
$ws.Range("I79").Value = 321
$ws.Range("I83").Value = 246
$ws.Range("I84").Value = 108
$ws.Range("I85").Value = 573
$ws.Range("I86").Value = 73
$ws.Range("I88").Value = 112
$ws.Range("I90").Value = 154
$ws.Range("I94").Value = 109
$ws.Range("I95").Value = 197
$ws.Range("I96").Value = 140
$ws.Range("I98").Value = 78
$ws.Range("I99").Value = 227
$ws.Range("G101").Value = 24657
$ws.Range("H101").Value = 25978
$ws.Range("I101").Value = 12347

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("I2").Value = 148
$ws.Range("I3").Value = 231
$ws.Range("I7").Value = 573

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("I3").Value = 94
$ws.Range("I6").Value = 63
$ws.Range("I7").Value = 270

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("I2").Value = 238
$ws.Range("I3").Value = 209
$ws.Range("I5").Value = 24
$ws.Range("I6").Value = 243
$ws.Range("I7").Value = 757

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("I6").Value = 62
$ws.Range("I7").Value = 138

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("G4").Value = 39
$ws.Range("G7").Value = 708

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("I3").Value = 42
$ws.Range("I5").Value = 1
$ws.Range("I6").Value = 53
$ws.Range("I7").Value = 140

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("I2").Value = 126
$ws.Range("I7").Value = 391

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("I5").Value = 6
$ws.Range("I7").Value = 227

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I2").Value = 115
$ws.Range("I3").Value = 169
$ws.Range("I7").Value = 476

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("I6").Value = 29
$ws.Range("I7").Value = 108

$ws = $wb.Worksheets.Item("New City")
$ws.Range("I2").Value = 89
$ws.Range("I3").Value = 75
$ws.Range("I6").Value = 85
$ws.Range("I7").Value = 273

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("I5").Value = 11
$ws.Range("I7").Value = 246

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("I3").Value = 73
$ws.Range("I7").Value = 197

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I2").Value = 131
$ws.Range("I3").Value = 199
$ws.Range("I7").Value = 557

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I6").Value = 142
$ws.Range("I7").Value = 276

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I2").Value = 240
$ws.Range("I3").Value = 277
$ws.Range("I4").Value = 39
$ws.Range("I6").Value = 216
$ws.Range("I7").Value = 802

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("I2").Value = 127
$ws.Range("I6").Value = 89
$ws.Range("I7").Value = 324

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("I3").Value = 24
$ws.Range("I7").Value = 92

$ws = $wb.Worksheets.Item("River North")
$ws.Range("I3").Value = 46
$ws.Range("I7").Value = 188

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("I2").Value = 37
$ws.Range("I7").Value = 85

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I3").Value = 147
$ws.Range("I7").Value = 426

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("I6").Value = 70
$ws.Range("I7").Value = 174

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("I3").Value = 37
$ws.Range("I6").Value = 47
$ws.Range("I7").Value = 140

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("I2").Value = 6
$ws.Range("I7").Value = 28

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("I3").Value = 105
$ws.Range("I6").Value = 96
$ws.Range("I7").Value = 321

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("I3").Value = 20
$ws.Range("I7").Value = 82

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("I6").Value = 62
$ws.Range("I7").Value = 109

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("I6").Value = 51
$ws.Range("I7").Value = 146

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("I4").Value = 5
$ws.Range("I6").Value = 48
$ws.Range("I7").Value = 78

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("I2").Value = 14
$ws.Range("I4").Value = 13
$ws.Range("I7").Value = 56

$ws = $wb.Worksheets.Item("East Village")
$ws.Range("I6").Value = 9
$ws.Range("I7").Value = 15

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("I6").Value = 35
$ws.Range("I7").Value = 112

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("I4").Value = 15
$ws.Range("I7").Value = 110

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("I4").Value = 39
$ws.Range("I7").Value = 73

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("I3").Value = 15
$ws.Range("I7").Value = 43

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("I3").Value = 34
$ws.Range("I7").Value = 154

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("I4").Value = 12
$ws.Range("I7").Value = 120

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("I3").Value = 11
$ws.Range("I7").Value = 34

$ws = $wb.Worksheets.Item("Beverly")
$ws.Range("I3").Value = 3
$ws.Range("I7").Value = 24
